$wb = $excel.ActiveWorkbook

# Rename sheets to the new "V_" prefixed naming convention
$wb.Worksheets.Item("HomePage").Name = "V_HomePage"
$wb.Worksheets.Item("DataModelPage").Name = "V_DataModelPage"
$wb.Worksheets.Item("ResourcesPage").Name = "V_ResourcesPage"
$wb.Worksheets.Item("AboutAboutPage").Name = "V_AboutAboutPage"

# Move the selection on the AboutAboutPage tab before activating it
$aboutSheet = $wb.Worksheets.Item("V_AboutAboutPage")
$aboutSheet.Activate()
$aboutSheet.Range("B33").Select()

# Activate the AboutAboutPage tab (4th tab, 0-indexed activeTab = 3)
$wb.Worksheets.Item("V_AboutAboutPage").Activate()
